$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new quarter "Q8" column (J) was added to the report.
# Copy the formatting from the existing last header cell (I1) onto the new header cell (J1)
# so the new column matches the look of the rest of the header row, then set its label.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").Value = "Q8"

# Data rows 2-16 were recomputed (the underlying naive-error calculation changed),
# and rows 4 and 8 now have enough trailing history to populate extra columns (G:J / G:I).
$ws.Range("B2").Value = 0.168549046501717
$ws.Range("C2").Value = -0.04441333468075517
$ws.Range("D2").Value = -0.7572115014490061
$ws.Range("E2").Value = -0.441590891698371
$ws.Range("F2").Value = 0.04279636911810641
$ws.Range("G2").Value = 0.06926039754389782
$ws.Range("H2").Value = -0.190909471160479
$ws.Range("I2:J2").ClearContents()
$ws.Range("B3").Value = -0.148691993218036
$ws.Range("C3").Value = -0.8614901599862869
$ws.Range("D3").Value = -0.5458695502356519
$ws.Range("E3").Value = -0.06148228941917439
$ws.Range("F3").Value = -0.03501826099338298
$ws.Range("G3").Value = -0.2951881296977598
$ws.Range("H3:J3").ClearContents()
$ws.Range("B4").Value = -0.7001939218680397
$ws.Range("C4").Value = -0.3845733121174047
$ws.Range("D4").Value = 0.09981394869907274
$ws.Range("E4").Value = 0.1262779771248642
$ws.Range("F4").Value = -0.1338918915795126
$ws.Range("G4").Value = 0.067006525927944
$ws.Range("H4").Value = -0.3988798524297023
$ws.Range("I4").Value = -0.006388367470709932
$ws.Range("J4").Value = -0.05594074298689113
$ws.Range("B5").Value = -0.0532934954405796
$ws.Range("C5").Value = 0.4310937653758978
$ws.Range("D5").Value = 0.4575577938016893
$ws.Range("E5").Value = 0.1973879250973124
$ws.Range("F5").Value = 0.3982863426047691
$ws.Range("G5").Value = -0.06760003575287715
$ws.Range("H5").Value = 0.3248914492061151
$ws.Range("I5").Value = 0.2753390736899339
$ws.Range("J5").ClearContents()
$ws.Range("B6").Value = 0.05953170440083436
$ws.Range("C6").Value = 0.08599573282662576
$ws.Range("D6").Value = -0.174174135877751
$ws.Range("E6").Value = 0.02672428162970562
$ws.Range("F6").Value = -0.4391620967279406
$ws.Range("G6").Value = -0.04667061176894832
$ws.Range("H6").Value = -0.09622298728512951
$ws.Range("I6:J6").ClearContents()
$ws.Range("B7").Value = 0.027557006744978
$ws.Range("C7").Value = -0.2326128619593988
$ws.Range("D7").Value = -0.03171444445194214
$ws.Range("E7").Value = -0.4976008228095884
$ws.Range("F7").Value = -0.1051093378505961
$ws.Range("G7").Value = -0.1546617133667773
$ws.Range("H7:J7").ClearContents()
$ws.Range("B8").Value = -0.185259264801978
$ws.Range("C8").Value = 0.01563915270547871
$ws.Range("D8").Value = -0.4502472256521675
$ws.Range("E8").Value = -0.05775574069317523
$ws.Range("F8").Value = -0.1073081162093564
$ws.Range("G8").Value = 0.2097631713009834
$ws.Range("H8").Value = 0.08508340067688215
$ws.Range("I8").Value = 0.04125054597109605
$ws.Range("J8").ClearContents()
$ws.Range("B9").Value = 0.1016561395759164
$ws.Range("C9").Value = -0.3642302387817298
$ws.Range("D9").Value = 0.02826124617726248
$ws.Range("E9").Value = -0.02129112933891871
$ws.Range("F9").Value = 0.2957801581714211
$ws.Range("G9").Value = 0.1711003875473198
$ws.Range("H9").Value = 0.1272675328415337
$ws.Range("I9:J9").ClearContents()
$ws.Range("B10").Value = -0.3543200805324755
$ws.Range("C10").Value = 0.0381714044265168
$ws.Range("D10").Value = -0.01138097108966439
$ws.Range("E10").Value = 0.3056903164206754
$ws.Range("F10").Value = 0.1810105457965742
$ws.Range("G10").Value = 0.1371776910907881
$ws.Range("H10:J10").ClearContents()
$ws.Range("B11").Value = 0.1247600422622814
$ws.Range("C11").Value = 0.07520766674610024
$ws.Range("D11").Value = 0.39227895425644
$ws.Range("E11").Value = 0.2675991836323388
$ws.Range("F11").Value = 0.2237663289265527
$ws.Range("G11:J11").ClearContents()
$ws.Range("B12").Value = 0.015958518250201
$ws.Range("C12").Value = 0.3330298057605408
$ws.Range("D12").Value = 0.2083500351364396
$ws.Range("E12").Value = 0.1645171804306535
$ws.Range("F12:J12").ClearContents()
$ws.Range("B13").Value = 0.2879850433121589
$ws.Range("C13").Value = 0.1633052726880577
$ws.Range("D13").Value = 0.1194724179822716
$ws.Range("E13:J13").ClearContents()
$ws.Range("B14").Value = 0.0315981439370237
$ws.Range("C14").Value = -0.0122347107687624
$ws.Range("D14:J14").ClearContents()
$ws.Range("B15").Value = 0.05023326629364772
$ws.Range("C15:J15").ClearContents()
$ws.Range("B16:J16").ClearContents()
